{"js": "// Fix the \"Count Real Numbers\" hints paragraph: the hint used to describe a\n// `SortedDictionary<double, int>` (leftover from the integer-counting task);\n// replace it with the right description for this (key=nums, value=count)\n// sorted dictionary used in the \"Count Real Numbers\" problem.\nconst body = context.document.body;\n\n// 1) \"SortedDictionary<double,\" -> \"sorted dicrtionary (key=nums\" (keep CodeChar style)\nconst part1 = body.search(\"SortedDictionary<double,\", { matchCase: true });\npart1.load(\"items\");\nawait context.sync();\nif (part1.items.length > 0) {\n  const hit1 = part1.items[0];\n  hit1.insertText(\"sorted dicrtionary (key=nums\", \"Replace\");\n  const comma = hit1.insertText(\",\", \"After\");\n  comma.style = \"CodeChar\";\n  await context.sync();\n}\n\n// 2) \"int>\" -> \"value=count) \" (keep CodeChar style, note trailing space\n//    moves into this run).\nconst part2 = body.search(\"int>\", { matchCase: true });\npart2.load(\"items\");\nawait context.sync();\nif (part2.items.length > 0) {\n  const hit2 = part2.items[0];\n  hit2.insertText(\"value=count) \", \"Replace\");\n  await context.sync();\n}\n\n// 2b) the old \" named \" (leading space + \"named\" + trailing space) loses its\n//     leading space (it moved into the \"value=count) \" run above) so it\n//     becomes \"named \" (keep the plain/no style formatting).\nconst part2b = body.search(\" named \", { matchCase: true });\npart2b.load(\"items\");\nawait context.sync();\nif (part2b.items.length > 0) {\n  const hit2b = part2b.items[0];\n  hit2b.insertText(\"named \", \"Replace\");\n  await context.sync();\n}\n\n// 3) \"counts.Keys\" -> \"counts.keys()\" (still CodeChar style), and move the\n//    \"_GoBack\" bookmark to sit right after it (it used to sit alone in the\n//    trailing empty paragraph at the end of the document).\nconst part3 = body.search(\"counts.Keys\", { matchCase: true });\npart3.load(\"items\");\nawait context.sync();\nif (part3.items.length > 0) {\n  const hit3 = part3.items[0];\n  hit3.insertText(\"counts.keys()\", \"Replace\");\n  await context.sync();\n\n  context.document.deleteBookmark(\"_GoBack\");\n  const afterHit3 = body.search(\"counts.keys()\", { matchCase: true });\n  afterHit3.load(\"items\");\n  await context.sync();\n  if (afterHit3.items.length > 0) {\n    const endRange = afterHit3.items[0].getRange(\"End\");\n    endRange.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Fix the \"Count Real Numbers\" hints paragraph: the hint used to describe a\n# `SortedDictionary<double, int>` (leftover from the integer-counting task);\n# replace it with the right description for this (key=nums, value=count)\n# sorted dictionary used in the \"Count Real Numbers\" problem.\n$d = $word.ActiveDocument\n\n# 1) \"SortedDictionary<double,\" -> \"sorted dicrtionary (key=nums\" (keep CodeChar style)\n$rng = $d.Content\n$found = $rng.Find.Execute(\"SortedDictionary<double,\")\nif ($found) {\n    $rng.Text = \"sorted dicrtionary (key=nums\"\n    $endPt = $rng.Duplicate\n    $endPt.Collapse(0)  # wdCollapseEnd\n    $endPt.InsertAfter(\",\")\n    $endPt.Style = \"CodeChar\"\n}\n\n# 2) \"int>\" -> \"value=count) \" (keep CodeChar style, note trailing space\n#    moves into this run).\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"int>\")\nif ($found2) {\n    $rng2.Text = \"value=count) \"\n}\n\n# 2b) the old \" named \" (leading space + \"named\" + trailing space) loses its\n#     leading space (it moved into the \"value=count) \" run above) so it\n#     becomes \"named \" (keep the plain/no style formatting).\n$rng2b = $d.Content\n$found2b = $rng2b.Find.Execute(\" named \")\nif ($found2b) {\n    $rng2b.Text = \"named \"\n}\n\n# 3) \"counts.Keys\" -> \"counts.keys()\" (still CodeChar style), and move the\n#    \"_GoBack\" bookmark to sit right after it (it used to sit alone in the\n#    trailing empty paragraph at the end of the document).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$rng3 = $d.Content\n$found3 = $rng3.Find.Execute(\"counts.Keys\")\nif ($found3) {\n    $rng3.Text = \"counts.keys()\"\n    $endPt3 = $rng3.Duplicate\n    $endPt3.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $endPt3)\n}\n"}
